$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $text) {
    # Prefix with an apostrophe so Excel stores the literal text instead
    # of auto-converting date-shaped strings ("2017-05-27", ...) into
    # serial date numbers. ClearFormats() afterwards drops the transient
    # quote-prefix formatting again so the cell falls back to the sheet's
    # default (unstyled) cell format, same as the rest of the table.
    $r = $ws.Range($range)
    $r.Value = "'" + $text
    $r.ClearFormats()
}

# Row 10: "# 65 LRPM"
Set-TextCell "A10" "# 65 LRPM"
Set-TextCell "B10" "Kyiv cert"
Set-TextCell "C10" "2017-05-27"
Set-TextCell "D10" "2017-05-31"
Set-TextCell "E10" ""
Set-TextCell "F10" ""
Set-TextCell "G10" "Igor  Gnes"
Set-TextCell "H10" "2017-05-22"

# Row 11: "# 66 UTLS"
Set-TextCell "A11" "# 66 UTLS"
Set-TextCell "B11" "Pquur60GiHFw/u8gsV1BcC+YkI6X5t+6yMJjORrA4RL4dLPXzFPnLVTHpCvxmdUXXFLpov3e+db2Kbod+cuwyRLt9HsUQ6uSTrhD2e1LObRi1SG47Y3ECeYtD2LNVogHybQyyjZwi8IoLPySWms1+/1NsHa7CmduFuutWfk6qNw="
Set-TextCell "C11" "2017-05-29"
Set-TextCell "D11" "2017-06-01"
Set-TextCell "E11" ""
Set-TextCell "F11" ""
Set-TextCell "G11" "Igor  Gnes"
Set-TextCell "H11" "2017-05-22"

# Re-run best-fit auto sizing on every column now that the new rows (in
# particular the long encoded value in column B) have been added.
$ws.Columns.Item(1).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(2).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(3).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(4).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(5).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(6).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(7).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(8).EntireColumn.AutoFit() | Out-Null
